# Adds a new "2022-Q1" sheet (fund-holdings detail) right before the "总计"
# (totals) sheet, and refreshes "总计" with a new leading row for 2022-Q1.
#
# Implementation note: rather than inserting a brand-new blank sheet and
# retyping the *existing* "总计" data into it, we repurpose the worksheet
# that is currently named "总计" (it already sits in the right slot, tab 6)
# as the new "2022-Q1" detail sheet, and create a brand-new worksheet named
# "总计" after it. That mirrors the sheetId bookkeeping in the target diff:
# the pre-existing sheet (sheetId 6) keeps its id under its new name
# "2022-Q1", and the freshly-created sheet gets the next id (7) under the
# name "总计".
#
# IMPORTANT engine quirk: worksheet object references handed out by this
# COM shim are resolved against the *live tab position*, not a stable
# identity. Any call that changes the sheet count or order (Worksheets.Add,
# Worksheet.Move, Worksheet.Delete) silently repoints every worksheet
# variable already in scope at whatever now sits at that old index. So we
# never hold a worksheet variable across such a call - we always re-fetch
# by name with Worksheets.Item(...) immediately before using it again.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: repurpose the current "总计" sheet (tab 6) into "2022-Q1" and
# lay down the fund-holdings table in it. No Add/Move/Delete happens
# anywhere in this block, so it's safe to reuse $q1 throughout.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(6)
$q1.Name = "2022-Q1"

# Wipe the old 日期/持有数量/持有市值 total-table content + formatting so we
# can lay down the fund-holdings table in its place.
$q1.Cells.Clear()

# Header row (row 1) for the fund-holdings table.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Fund rows: code, name, scale, total stock position, position pct,
# held value (亿元), position rank.
$fundRows = @(
    @("005739", "富国转型机遇混合", "39.36", "80.18", "3.24", "1.2753", 5),
    @("100056", "富国低碳环保混合", "27.82", "82.05", "3.27", "0.9097", 6),
    @("506005", "博时科创板三年定期开放混合", "22.84", "96.44", "3.91", "0.8930", 5),
    @("506003", "富国科创板两年定期开放混合", "23.60", "98.48", "2.53", "0.5971", 5),
    @("012888", "工银兴瑞一年持有期混合A", "22.32", "76.39", "2.56", "0.5714", 10),
    @("213003", "宝盈策略增长混合", "10.28", "94.38", "4.73", "0.4862", 8),
    @("011212", "富国稳健策略6个月持有期混合A", "12.61", "83.04", "3.42", "0.4313", 6),
    @("006527", "富国优质发展混合A", "13.71", "79.35", "2.91", "0.3990", 7),
    @("005535", "泰信竞争优选灵活配置混合", "10.07", "90.21", "3.83", "0.3857", 5),
    @("213002", "宝盈泛沿海增长混合", "5.10", "93.76", "6.71", "0.3422", 4),
    @("290006", "泰信蓝筹精选混合", "7.62", "88.52", "3.29", "0.2507", 9),
    @("006528", "富国优质发展混合C", "4.04", "79.35", "2.91", "0.1176", 7),
    @("012358", "汇丰晋信医疗先锋混合型证券投资基金A", "2.12", "85.40", "3.31", "0.0702", 9),
    @("012889", "工银兴瑞一年持有期混合C", "2.22", "76.39", "2.56", "0.0568", 10),
    @("011213", "富国稳健策略6个月持有期混合C", "1.14", "83.04", "3.42", "0.0390", 6),
    @("010756", "兴华永兴混合A", "0.35", "94.57", "4.48", "0.0157", 6),
    @("003749", "创金合信鑫收益灵活配置混合A", "0.65", "51.22", "1.17", "0.0076", 9),
    @("006906", "创金合信鑫收益灵活配置混合E", "0.65", "51.22", "1.17", "0.0076", 9),
    @("012359", "汇丰晋信医疗先锋混合型证券投资基金C", "0.11", "85.40", "3.31", "0.0036", 9),
    @("010757", "兴华永兴混合C", "0.01", "94.57", "4.48", "0.0004", 6),
    @("003750", "创金合信鑫收益灵活配置混合C", "0.02", "51.22", "1.17", "0.0002", 9)
)

for ($i = 0; $i -lt $fundRows.Count; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $q1.Cells.Item($r, 1).Value = $i

    $q1.Cells.Item($r, 2).Value = "'" + $row[0]
    $q1.Cells.Item($r, 2).Style = "Normal"

    $q1.Cells.Item($r, 3).Value = $row[1]

    $q1.Cells.Item($r, 4).Value = "'" + $row[2]
    $q1.Cells.Item($r, 4).Style = "Normal"

    $q1.Cells.Item($r, 5).Value = "'" + $row[3]
    $q1.Cells.Item($r, 5).Style = "Normal"

    $q1.Cells.Item($r, 6).Value = "'" + $row[4]
    $q1.Cells.Item($r, 6).Style = "Normal"

    $q1.Cells.Item($r, 7).Value = "'" + $row[5]
    $q1.Cells.Item($r, 7).Style = "Normal"

    $q1.Cells.Item($r, 8).Value = $row[6]
}

# Formatting - header row + index column match the bold / bordered /
# centered look used by every other quarter tab. Pull that format
# straight from the "2021-Q4" tab, which has the identical 8-column
# fund-holdings layout. (Copy/PasteSpecial don't reorder sheets, so $q1
# stays valid through this.)
$fmtSource = $wb.Worksheets.Item("2021-Q4")
$fmtSource.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$fmtSource.Range("A2").Copy()
$q1.Range("A2:A22").PasteSpecial(-4122)

$q1.Range("A1").Select()

# ---------------------------------------------------------------------
# Step 2: append a fresh "总计" sheet after "2022-Q1" with the running
# per-quarter totals table, now including the new 2022-Q1 row on top.
#
# Worksheets.Add()/Move() reorder tabs, which invalidates every
# worksheet variable obtained earlier (including $q1, $fmtSource) - so
# from this point on we re-fetch everything we need by name right
# before we use it, and we don't touch $q1 / $fmtSource again.
# ---------------------------------------------------------------------
$newTotal = $wb.Worksheets.Add()
$newTotal.Name = "总计"

# Re-fetch "2022-Q1" by name since Worksheets.Add() just shifted tab
# indices - any reference taken before this point may now resolve to a
# different sheet.
$q1Again = $wb.Worksheets.Item("2022-Q1")
$newTotal.Move([System.Reflection.Missing]::Value, $q1Again)

# Move() also reshuffles tab order, so re-fetch "总计" fresh by name
# before writing anything into it.
$total = $wb.Worksheets.Item("总计")

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @("2022-Q1", 21, 6.86),
    @("2021-Q4", 16, 4.67),
    @("2021-Q3", 21, 4.79),
    @("2021-Q2", 63, 16.6),
    @("2021-Q1", 23, 6.68),
    @("2020-Q4", 7, 1.61)
)

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = "'" + $row[0]
    $total.Cells.Item($r, 2).Style = "Normal"
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
}

$fmtSource2 = $wb.Worksheets.Item("2021-Q4")
$fmtSource2.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$fmtSource2.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Range("A1").Select()
